$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "PYTHON"

$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Feuil1"
